$wb = $excel.ActiveWorkbook

$values = @{
    "forprofit"  = @(-0.03466775208093566, -0.01302760197228154, -0.07338944145858986, 0.002206286492960288, 0.1201470517820603)
    "nonprofit"  = @(-0.003085812737391978, 0.01374037238034471, -0.0468626910563708, 0.001042362587145307, 0.03537888105551978)
    "government" = @(-0.001072047097747203, 0.01050064234067945, -0.01826618174083947, 0.005880029066311169, 0.003186633733128033)
    "total"      = @(-0.02148302485726171, -0.0008693896351791266, -0.06070361124259566, 0.005380378932626438, 0.07798163046410518)
}

foreach ($sheetName in $values.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowValues = $values[$sheetName]
    for ($i = 0; $i -lt $rowValues.Length; $i++) {
        $rowNum = $i + 2
        $ws.Cells.Item($rowNum, 2).Value = $rowValues[$i]
    }
}
